$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the "Espinaca" series at
# Feria Lagunitas de Puerto Montt. It belongs chronologically right after
# the existing row 20 (date 2022-09-29 / serial 44838), so insert a new
# row at position 21, which pushes the old rows 21-55 down to 22-56.
$ws.Rows("21:21").Insert()

# Fill in the newly inserted row 21 with the new observation's data.
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C21").Value = "Los Lagos"
$ws.Range("D21").Value = 45002
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 100112012
$ws.Range("G21").Value = "Espinaca"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 25
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 15000
$ws.Range("N21").Value = "$/cuna 10 kilos"
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 1500
$ws.Range("Q21").Value = 10
$ws.Range("R21").Value = "Hortaliza"
